$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remove the "Dim Name Here!" label that had been placed in L5 ---
# (this also drops the now-unused shared string from sst on save)
$ws.Range("L5").ClearContents()

# --- Widen columns C and D to fit their (now longer) header/content ---
$ws.Columns.Item(3).ColumnWidth = 23.736979166666668
$ws.Columns.Item(4).ColumnWidth = 34.451822916666664

# --- Update the sheet view: scroll so column E is the leftmost visible
#     column, and move the selection to L5 (where the edit happened) ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("L5").Select()
